# Fix developer guide references to address book
# (UndoRedoActivityDiagram: rename "address book" domain example to
#  "expense tracker" in the two activity-diagram annotation shapes.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 7 "TextBox 47": "[command commits address book]"
#   -> "[command commits expense tracker]"
$shpGuard = $s.Shapes.Item(7)
$guardText = $shpGuard.TextFrame.TextRange
$guardRun = $guardText.Characters(2, 30)
$guardRun.Text = "command commits expense tracker]"

# Shape 8 "Rectangle: Rounded Corners 50":
#   "Purge redundant states and then save address book to addressBookStateList "
#   -> "Purge redundant states and then save expense tracker to expenseTrackerStateList "
$shpNote = $s.Shapes.Item(8)
$noteText = $shpNote.TextFrame.TextRange

# Edit right-to-left so earlier character offsets stay valid.
$stateListRun = $noteText.Characters(54, 20)
$stateListRun.Text = "expenseTrackerStateList"

$saveRun = $noteText.Characters(1, 53)
$saveRun.Text = "Purge redundant states and then save expense tracker to "
